$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds plain numeric-looking text (e.g. "37.868.59", "0.100").
# Excel auto-converts plain decimal-looking input typed into a General-format
# cell into a real number (losing formatting like trailing zeros), so force the
# cell to Text before assigning, then clear formats to restore the original
# (unstyled / General) appearance while keeping the value as literal text.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.868.59"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +1.54%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.088.13"
$ws.Range("D3").ClearFormats()
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "232.86"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.45%  "
$ws.Range("E6").Value = "  +0.23%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "57.39"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +1.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.388"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +1.85%  "
$ws.Range("E10").Value = "  +2.23%  "
$ws.Range("E11").Value = "  +2.93%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.382.17"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.58%  "
$ws.Range("E13").Value = "  -1.39%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.10"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +2.31%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.762"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.65%  "
$ws.Range("E16").Value = "  +2.35%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.077.21"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.66%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.785.61"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.37%  "
$ws.Range("E19").Value = "  -1.90%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.76"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.84%  "
$ws.Range("E21").Value = "  +1.22%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "228.36"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.00%  "
$ws.Range("E23").Value = "  -0.09%  "
$ws.Range("E24").Value = "  -1.70%  "
$ws.Range("E25").Value = "  -0.21%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "170.95"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +2.24%  "
$ws.Range("E27").Value = "  +9.70%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.96"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +2.42%  "
$ws.Range("E29").Value = "  +0.51%  "
$ws.Range("E30").Value = "  +2.30%  "
$ws.Range("E31").Value = "  +1.44%  "
$ws.Range("E32").Value = "  +4.01%  "
$ws.Range("E33").Value = "  +1.84%  "
$ws.Range("E34").Value = "  +1.19%  "
$ws.Range("E35").Value = "  +0.46%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.84"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +3.86%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.39"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +4.92%  "
$ws.Range("E38").Value = "  -0.26%  "
$ws.Range("E39").Value = "  -4.08%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.100"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +7.05%  "
$ws.Range("E41").Value = "  -0.89%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "97.06"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +1.07%  "
$ws.Range("E43").Value = "  +0.83%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.450.95"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -1.16%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.17"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.06%  "
$ws.Range("E46").Value = "  +3.41%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "15.72"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +4.87%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.06"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -6.55%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.39"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +3.67%  "
$ws.Range("E50").Value = "  +1.56%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.277.37"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.90%  "
